# Add support for HPC Cache: insert a new "hpc_cache" worksheet between
# "blob" and "anf", populate it with the HPC Cache cost table, and make it
# the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Find the "anf" sheet so we can insert the new sheet right before it.
$anfSheet = $wb.Worksheets.Item("anf")

$newSheet = $wb.Worksheets.Add($anfSheet)
$newSheet.Name = "hpc_cache"

$newSheet.Range("A1").Value = "HPC Cache"

$newSheet.Range("A5").Value = "Throughput_GBps"
$newSheet.Range("B5").Value = "Capacity_small_TiB"
$newSheet.Range("C5").Value = "Capacity_medium_TiB"
$newSheet.Range("D5").Value = "Capacity_large_TiB"
$newSheet.Range("E5").Value = "cost_small_per_month"
$newSheet.Range("F5").Value = "cost_medium_per_month"
$newSheet.Range("G5").Value = "cost_large_per_month"

$newSheet.Range("A3").Value = "cost per Month (PAYGO)"

$newSheet.Range("A6").Value = 2
$newSheet.Range("B6").Value = 3
$newSheet.Range("C6").Value = 6
$newSheet.Range("D6").Value = 12
$newSheet.Range("E6").Value = 4881.82
$newSheet.Range("F6").Value = 5581.5
$newSheet.Range("G6").Value = 6980.86

$newSheet.Range("A7").Value = 4
$newSheet.Range("B7").Value = 6
$newSheet.Range("C7").Value = 12
$newSheet.Range("D7").Value = 24
$newSheet.Range("E7").Value = 9763.64
$newSheet.Range("F7").Value = 11163
$newSheet.Range("G7").Value = 13961.71

$newSheet.Range("A8").Value = 8
$newSheet.Range("B8").Value = 12
$newSheet.Range("C8").Value = 24
$newSheet.Range("D8").Value = 48
$newSheet.Range("E8").Value = 19527.28
$newSheet.Range("F8").Value = 22325.99
$newSheet.Range("G8").Value = 27923.32

$newSheet.Columns.Item(1).ColumnWidth = 15.59765625
$newSheet.Columns.Item(2).ColumnWidth = 16.796875
$newSheet.Columns.Item(3).ColumnWidth = 19.33203125
$newSheet.Columns.Item(4).ColumnWidth = 17.46484375
$newSheet.Columns.Item(5).ColumnWidth = 19.46484375
$newSheet.Columns.Item(6).ColumnWidth = 22.265625
$newSheet.Columns.Item(7).ColumnWidth = 20.6640625

$newSheet.Activate()
$newSheet.Range("C17").Select()
$excel.ActiveWindow.ScrollColumn = 2
